$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated Hours for Mike (column D, Week 1) and Patrick (column F, Week 1)
$ws.Range("D4").Value = 2
$ws.Range("F4").Value = 2

# Recalculate so dependent formulas (E4, G4, E5:E13, etc.) update their cached values
$excel.Calculate()

# Move the active selection to H4, as in the saved workbook
$ws.Range("H4").Select()
